# Updated bento multifilter test scripts
# The "startup" sheet (ActiveSheet) previously listed three tabs (CasesTab,
# SamplesTab, FilesTab) each pointing at TC01 Akita-breed-filter queries.
# It is trimmed down to a single CasesTab row, and the referenced
# Neo4j/Web data file names are switched over to the new TC03 multi-filter
# (Study/Breed/Sex) test file names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the SamplesTab (row 3) and FilesTab (row 4) rows entirely,
# shifting everything below them up (rows 5-13 become rows 3-11).
$ws.Rows("3:4").Delete()

# Point the remaining CasesTab row at the new TC03 manifest file names.
$ws.Range("D2").Value = "TC03_Canine_E2E_MultipleFilters-Study_Breed_Sex_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC03_Canine_E2E_MultipleFilters-Study_Breed_Sex_WebData.xlsx"

# Update the saved selection/view state to match the edited sheet.
$ws.Range("C7").Select()
